$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 148; this shifts the existing rows 148-198
# down to 149-199 (and grows the used range to A1:R199).
$ws.Rows(148).Insert()

# Populate the newly inserted row 148 with the new price record
# (Papa / Patagonia / "1a (guarda)") that was added between the prior
# "1a (cosecha)" entry and the rest of the existing data.
$ws.Range("A148").Value = 7
$ws.Range("B148").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C148").Value = "Ñuble"
$ws.Range("D148").Value = 44468
$ws.Range("E148").Value = 16
$ws.Range("F148").Value = 100114001
$ws.Range("G148").Value = "Papa"
$ws.Range("H148").Value = "Patagonia"
$ws.Range("I148").Value = "1a (guarda)"
$ws.Range("J148").Value = 300
$ws.Range("K148").Value = 7000
$ws.Range("L148").Value = 7500
$ws.Range("M148").Value = 7250
$ws.Range("N148").Value = "$/saco 25 kilos"
$ws.Range("O148").Value = "Provincia de Diguillín"
$ws.Range("P148").Value = 290
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = "Hortaliza"
